# Updates cryptos list values (prices / 1h volume %) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.723.91"
Set-TextValue "E2" "  +2.88%  "
Set-TextValue "D3" "3.127.64"
Set-TextValue "E3" "  +1.62%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "588.85"
Set-TextValue "E5" "  +1.58%  "
Set-TextValue "D6" "147.38"
Set-TextValue "E6" "  +3.63%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "3.124.21"
Set-TextValue "E8" "  +1.81%  "
Set-TextValue "E9" "  +0.71%  "
Set-TextValue "E10" "  +14.58%  "
Set-TextValue "E11" "  -0.14%  "
Set-TextValue "E12" "  +0.79%  "
Set-TextValue "D13" "0.0000252"
Set-TextValue "E13" "  +4.65%  "
Set-TextValue "D14" "37.23"
Set-TextValue "E14" "  +5.77%  "
Set-TextValue "E15" "  -0.72%  "
Set-TextValue "D16" "3.644.50"
Set-TextValue "E16" "  +1.62%  "
Set-TextValue "E17" "  -1.28%  "
Set-TextValue "D18" "63.622.23"
Set-TextValue "E18" "  +2.85%  "
Set-TextValue "D19" "3.125.72"
Set-TextValue "E19" "  +1.79%  "
Set-TextValue "D20" "465.71"
Set-TextValue "E20" "  +3.96%  "
Set-TextValue "E21" "  +2.49%  "
Set-TextValue "D22" "0.733"
Set-TextValue "E22" "  -0.25%  "
Set-TextValue "D23" "7.55"
Set-TextValue "E23" "  +1.44%  "
Set-TextValue "D24" "13.32"
Set-TextValue "E24" "  -3.38%  "
Set-TextValue "D25" "82.33"
Set-TextValue "E25" "  +0.75%  "
Set-TextValue "E26" "  +0.01%  "
Set-TextValue "D27" "8.99"
Set-TextValue "E27" "  +9.69%  "
Set-TextValue "E28" "  +1.91%  "
Set-TextValue "E29" "  -1.15%  "
Set-TextValue "E30" "  +0.00%  "
Set-TextValue "E31" "  +0.96%  "
Set-TextValue "D32" "27.13"
Set-TextValue "E32" "  +1.23%  "
Set-TextValue "E33" "  -3.95%  "
Set-TextValue "D34" "0.0₃0877"
Set-TextValue "E34" "  +10.41%  "
Set-TextValue "E35" "  +8.00%  "
Set-TextValue "B36" "dogwifhat"
Set-TextValue "C36" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D36" "3.45"
Set-TextValue "E36" "  +16.78%  "
Set-TextValue "B37" "Mantle"
Set-TextValue "C37" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D37" "1.05"
Set-TextValue "E37" "  +1.49%  "
Set-TextValue "E38" "  +1.27%  "
Set-TextValue "D39" "457.02"
Set-TextValue "E39" "  +8.91%  "
Set-TextValue "D40" "50.98"
Set-TextValue "E40" "  +1.65%  "
Set-TextValue "D41" "8.74"
Set-TextValue "E41" "  -0.66%  "
Set-TextValue "E42" "  +1.40%  "
Set-TextValue "D43" "2.912.00"
Set-TextValue "E43" "  -0.70%  "
Set-TextValue "D44" "0.278"
Set-TextValue "E44" "  +0.49%  "
Set-TextValue "E45" "  +2.27%  "
Set-TextValue "E46" "  +2.69%  "
Set-TextValue "D47" "127.55"
Set-TextValue "E47" "  +3.26%  "
Set-TextValue "E48" "  +2.15%  "
Set-TextValue "E50" "  +0.42%  "
Set-TextValue "D51" "24.73"
Set-TextValue "E51" "  +1.21%  "
